$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1474.7778
$ws.Range("I39").Value = 45.5
$ws.Range("J39").Value = 4333.3335
$ws.Range("K39").Value = 136.5
$ws.Range("L39").Value = 13000.0005
$ws.Range("M39").Value = 159.5
$ws.Range("N39").Value = -13592.0005
$ws.Range("H69").Value = 3583
$ws.Range("I69").Value = 3583
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 10749
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -9875
$ws.Range("H72").Value = 3583
$ws.Range("I72").Value = 3583
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 32247
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -27879
$ws.Range("H137").Value = 1722.0834
$ws.Range("I137").Value = 1838
$ws.Range("J137").Value = 1606.1666
$ws.Range("K137").Value = 5514
$ws.Range("L137").Value = 4818.4998
$ws.Range("M137").Value = -2964
$ws.Range("N137").Value = -9918.4998

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7800.2
$ws.Range("I32").Value = 7800.2
$ws.Range("K32").Value = 7800.2
$ws.Range("M32").Value = -7513.2
$ws.Range("H38").Value = 4500
$ws.Range("I38").Value = 4500
$ws.Range("K38").Value = 4500
$ws.Range("M38").Value = -4033
$ws.Range("H45").Value = 1333.1666
$ws.Range("I45").Value = 1199.8
$ws.Range("K45").Value = 1199.8
$ws.Range("M45").Value = -822.8
$ws.Range("H61").Value = 3441.4614
$ws.Range("I61").Value = 3212.7273
$ws.Range("K61").Value = 3212.7273
$ws.Range("M61").Value = -3000.7273
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H136").Value = 3441.4614
$ws.Range("I136").Value = 3212.7273
$ws.Range("K136").Value = 9638.1819
$ws.Range("M136").Value = -7088.1819

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2349.8333
$ws.Range("I105").Value = 1979.8
$ws.Range("J105").Value = 4200
$ws.Range("K105").Value = 1979.8
$ws.Range("L105").Value = 4200
$ws.Range("M105").Value = -232.8
$ws.Range("N105").Value = -7694

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 6816.6665
$ws.Range("I35").Value = 225
$ws.Range("J35").Value = 20000
$ws.Range("K35").Value = 225
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = 69
$ws.Range("N35").Value = -20588
$ws.Range("H38").Value = 4425
$ws.Range("I38").Value = 3850
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 3850
$ws.Range("L38").Value = 5000
$ws.Range("M38").Value = -3473
$ws.Range("N38").Value = -5754
$ws.Range("H46").Value = 4425
$ws.Range("I46").Value = 3850
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 3850
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -3639
$ws.Range("N46").Value = -5422
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H86").Value = 4749.75
$ws.Range("I86").Value = 4833
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 4833
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -3710
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 4749.75
$ws.Range("I89").Value = 4833
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 24165
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -18549
$ws.Range("N89").Value = -33732
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 4999
$ws.Range("K132").Value = 14997
$ws.Range("M132").Value = -12467
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 210.2
$ws.Range("J2").Value = 201
$ws.Range("L2").Value = 1206
$ws.Range("N2").Value = -1432

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 23.833334
$ws.Range("I2").Value = 23.833334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 23.833334
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 89.16666599999999
$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 5500
$ws.Range("K70").Value = 5500
$ws.Range("M70").Value = -5230
$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 5500
$ws.Range("K73").Value = 5500
$ws.Range("M73").Value = -4564
$ws.Range("H126").Value = 5055.5
$ws.Range("I126").Value = 5752.3335
$ws.Range("J126").Value = 2965
$ws.Range("K126").Value = 17257.0005
$ws.Range("L126").Value = 8895
$ws.Range("M126").Value = -14787.0005
$ws.Range("N126").Value = -13835

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1346.6666
$ws.Range("I32").Value = 1346.6666
$ws.Range("K32").Value = 1346.6666
$ws.Range("M32").Value = -1029.6666
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H136").Value = 1812.875
$ws.Range("I136").Value = 1167.1666
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 3501.4998
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -951.4998000000001
$ws.Range("N136").Value = -16350

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6856.7144
$ws.Range("I122").Value = 6599.4
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 19798.2
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -17348.2
$ws.Range("N122").Value = -27400
$ws.Range("H132").Value = 4164.4375
$ws.Range("I132").Value = 3283.1538
$ws.Range("K132").Value = 9849.4614
$ws.Range("M132").Value = -7319.4614
